$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "Comarca" column (D) metadata rows need to be reclassified from a
# measure to a dimension, matching the other reference-area columns.
$ws.Range("D3").Value = "sdmx-dimension:refArea"
$ws.Range("D4").Value = "dim"
$ws.Range("D5").Value = "URI-comarca"
